$d = $word.ActiveDocument

# The bullet about FourSquare currently reads "...filtered by Japanese
# restaurants...". It needs to become "...filtered by African restaurants...".
# The final XML splits the surrounding run into three pieces (with the
# document's _GoBack bookmark relocated to sit right after the new word),
# so first mark the two split points with bookmarks, then perform the text
# replacement, then drop the now-unneeded helper bookmark.

# Split point #1: between "venu" and "es and then filtered by Japanese".
$rng1 = $d.Content
$rng1.Find.Execute("venues and then filtered by Japanese", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint1 = $rng1.Start + 4
$d.Bookmarks.Add("_TempSplit", $d.Range($splitPoint1, $splitPoint1))

# Split point #2 (this is also where _GoBack should end up): right after
# "Japanese" (soon to be "African").
$rng2 = $d.Content
$rng2.Find.Execute("Japanese", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$splitPoint2 = $rng2.End
$d.Bookmarks.Add("_GoBack", $d.Range($splitPoint2, $splitPoint2))

# Replace the word itself. Adding the bookmarks above means this edit is
# confined to its own run, between the two split points.
$d.Content.Find.Execute("Japanese", $false, $false, $false, $false, $false, $true, 1, $false, "African", 2)

# Drop the helper bookmark now that it has done its job of splitting the run;
# the _GoBack bookmark (moved here from the end of the document) remains.
$d.Bookmarks("_TempSplit").Delete()
